$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G3").Value = "2016-09-03 20:47:09"
$wsZhCn.Range("H3").Value = "2016-09-03 20:47:00"
$wsZhCn.Range("K3").Value = "2016-09-03 20:47:30"
$wsDeDe.Range("H3").Value = "2016-09-03 20:47:09"
$wsDeDe.Range("K3").Value = "2016-09-03 20:47:37"
